$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-12 Thursday" "2025-06-13 Friday"

Replace-Text "304÷9=33, 7" "477÷8=59, 5"
Replace-Text "173÷7=24, 5" "629÷7=89, 6"
Replace-Text "779÷5=155, 4" "491÷5=98, 1"
Replace-Text "382÷6=63, 4" "826÷8=103, 2"
Replace-Text "882÷4=220, 2" "855÷8=106, 7"

Replace-Text "947÷5=189, 2" "775÷3=258, 1"
Replace-Text "391÷4=97, 3" "454÷9=50, 4"
Replace-Text "809÷7=115, 4" "842÷8=105, 2"
Replace-Text "530÷5=106, 0" "865÷7=123, 4"
Replace-Text "558÷6=93, 0" "664÷9=73, 7"

Replace-Text "106÷4=26, 2" "244÷3=81, 1"
Replace-Text "643÷7=91, 6" "717÷7=102, 3"
Replace-Text "805÷6=134, 1" "699÷6=116, 3"
Replace-Text "732÷5=146, 2" "918÷9=102, 0"
Replace-Text "876÷6=146, 0" "450÷9=50, 0"

Replace-Text "421÷6=70, 1" "766÷3=255, 1"
Replace-Text "414÷8=51, 6" "916÷4=229, 0"
Replace-Text "681÷3=227, 0" "986÷7=140, 6"
Replace-Text "649÷8=81, 1" "136÷5=27, 1"
Replace-Text "127÷2=63, 1" "206÷9=22, 8"

Replace-Text "855÷9=95, 0" "628÷9=69, 7"
Replace-Text "746÷4=186, 2" "713÷7=101, 6"
Replace-Text "658÷2=329, 0" "296÷7=42, 2"
Replace-Text "416÷2=208, 0" "989÷8=123, 5"
Replace-Text "453÷5=90, 3" "882÷9=98, 0"
